$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 29.5
$ws.Range("I6").Value = 29.5
$ws.Range("K6").Value = 88.5
$ws.Range("M6").Value = 23.5

$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H33").Value = 8066906
$ws.Range("I33").Value = 12500894
$ws.Range("K33").Value = 12500894
$ws.Range("M33").Value = -12500665

$ws.Range("H53").Value = 7030.737
$ws.Range("I53").Value = 1564
$ws.Range("J53").Value = 11950.8
$ws.Range("K53").Value = 1564
$ws.Range("L53").Value = 11950.8
$ws.Range("M53").Value = -927
$ws.Range("N53").Value = -13224.8

$ws.Range("H96").Value = 1579.2174
$ws.Range("I96").Value = 1377.6111
$ws.Range("J96").Value = 2305
$ws.Range("K96").Value = 4132.8333
$ws.Range("L96").Value = 6915
$ws.Range("M96").Value = -2759.8333
$ws.Range("N96").Value = -9661

$ws.Range("H99").Value = 830.5454999999999
$ws.Range("I99").Value = 897.6
$ws.Range("K99").Value = 2692.8
$ws.Range("M99").Value = -1194.8

$ws.Range("H107").Value = 1915.3334
$ws.Range("J107").Value = 2500
$ws.Range("L107").Value = 2500
$ws.Range("N107").Value = -6340

$ws.Range("H113").Value = 2364.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2364.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2364.5
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -8872.5

$ws.Range("H138").Value = 2421.5095
$ws.Range("I138").Value = 3350.7334
$ws.Range("J138").Value = 2054.7104
$ws.Range("K138").Value = 10052.2002
$ws.Range("L138").Value = 6164.1312
$ws.Range("M138").Value = -4912.200199999999
$ws.Range("N138").Value = -16444.1312

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 80000
$ws.Range("J7").Value = 80000
$ws.Range("L7").Value = 80000
$ws.Range("N7").Value = -80228

$ws.Range("H35").Value = 1699.6
$ws.Range("I35").Value = 874.75
$ws.Range("J35").Value = 4999
$ws.Range("K35").Value = 874.75
$ws.Range("L35").Value = 4999
$ws.Range("M35").Value = -468.75
$ws.Range("N35").Value = -5811

$ws.Range("H46").Value = 6384.2856
$ws.Range("I46").Value = 5912.6665
$ws.Range("K46").Value = 5912.6665
$ws.Range("M46").Value = -5593.6665

$ws.Range("H74").Value = 21787.586
$ws.Range("I74").Value = 1609.091
$ws.Range("K74").Value = 1609.091
$ws.Range("M74").Value = -735.0909999999999

$ws.Range("H77").Value = 21787.586
$ws.Range("I77").Value = 1609.091
$ws.Range("K77").Value = 8045.455
$ws.Range("M77").Value = -3677.455

$ws.Range("H133").Value = 60265
$ws.Range("J133").Value = 60265
$ws.Range("L133").Value = 60265
$ws.Range("N133").Value = -65325

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 50000004
$ws.Range("J19").Value = 100000000
$ws.Range("L19").Value = 100000000
$ws.Range("N19").Value = -100000346

$ws.Range("H134").Value = 54185.125
$ws.Range("I134").Value = 80147.30499999999
$ws.Range("K134").Value = 240441.915
$ws.Range("M134").Value = -237906.915

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 15913.546
$ws.Range("I31").Value = 981.36365
$ws.Range("K31").Value = 981.36365
$ws.Range("M31").Value = -686.36365

$ws.Range("H34").Value = 15913.546
$ws.Range("I34").Value = 981.36365
$ws.Range("K34").Value = 981.36365
$ws.Range("M34").Value = -779.36365

$ws.Range("H104").Value = 20000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 20000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 20000
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -25242

$ws.Range("H107").Value = 660.8823
$ws.Range("I107").Value = 766.26086
$ws.Range("K107").Value = 766.26086
$ws.Range("M107").Value = 1153.73914

$ws.Range("H109").Value = 18000
$ws.Range("J109").Value = 18000
$ws.Range("L109").Value = 18000
$ws.Range("N109").Value = -20080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3585611.8
$ws.Range("J4").Value = 400435.2
$ws.Range("L4").Value = 1201305.6
$ws.Range("N4").Value = -1201529.6

$ws.Range("H7").Value = 160.4
$ws.Range("J7").Value = 384.5
$ws.Range("L7").Value = 1153.5
$ws.Range("N7").Value = -1377.5

$ws.Range("H23").Value = 484
$ws.Range("I23").Value = 288
$ws.Range("J23").Value = 582
$ws.Range("K23").Value = 864
$ws.Range("L23").Value = 1746
$ws.Range("M23").Value = -629
$ws.Range("N23").Value = -2216

$ws.Range("H31").Value = 8499.5

$ws.Range("H80").Value = 20475
$ws.Range("J80").Value = 25712.5
$ws.Range("L80").Value = 77137.5
$ws.Range("N80").Value = -79009.5

$ws.Range("H83").Value = 20475
$ws.Range("J83").Value = 25712.5
$ws.Range("L83").Value = 231412.5
$ws.Range("N83").Value = -240772.5

$ws.Range("H92").Value = 278.45456
$ws.Range("I92").Value = 238.8
$ws.Range("J92").Value = 311.5
$ws.Range("K92").Value = 716.4000000000001
$ws.Range("L92").Value = 934.5
$ws.Range("M92").Value = 531.5999999999999
$ws.Range("N92").Value = -3430.5

$ws.Range("H107").Value = 938.4286
$ws.Range("I107").Value = 728.5714
$ws.Range("J107").Value = 1148.2858
$ws.Range("K107").Value = 2185.7142
$ws.Range("L107").Value = 3444.8574
$ws.Range("M107").Value = -265.7142000000003
$ws.Range("N107").Value = -7284.857400000001

$ws.Range("H121").Value = 431.5
$ws.Range("I121").Value = 377.8
$ws.Range("K121").Value = 1133.4
$ws.Range("M121").Value = 176.5999999999999

$ws.Range("H129").Value = 6062493.5
$ws.Range("I129").Value = 1878.8572
$ws.Range("J129").Value = 11365531
$ws.Range("K129").Value = 5636.571599999999
$ws.Range("L129").Value = 34096593
$ws.Range("M129").Value = -636.5715999999993
$ws.Range("N129").Value = -34106593

$ws.Range("H139").Value = 9519.333000000001
$ws.Range("I139").Value = 11921.929
$ws.Range("J139").Value = 4714.143
$ws.Range("K139").Value = 35765.787
$ws.Range("L139").Value = 14142.429
$ws.Range("M139").Value = -30625.787
$ws.Range("N139").Value = -24422.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2861.9
$ws.Range("I43").Value = 2123
$ws.Range("J43").Value = 3178.5715
$ws.Range("K43").Value = 2123
$ws.Range("L43").Value = 3178.5715
$ws.Range("M43").Value = -1972
$ws.Range("N43").Value = -3480.5715

$ws.Range("H46").Value = 1050
$ws.Range("I46").Value = 1050
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1050
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -894
$ws.Range("N46").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1990436.2
$ws.Range("I7").Value = 2651415
$ws.Range("J7").Value = 7500
$ws.Range("K7").Value = 2651415
$ws.Range("L7").Value = 7500
$ws.Range("M7").Value = -2651303
$ws.Range("N7").Value = -7724

$ws.Range("H16").Value = 100003490
$ws.Range("I16").Value = 200003940
$ws.Range("J16").Value = 3039.2
$ws.Range("K16").Value = 200003940
$ws.Range("L16").Value = 3039.2
$ws.Range("M16").Value = -200003770
$ws.Range("N16").Value = -3379.2

$ws.Range("H93").Value = 50006496
$ws.Range("I93").Value = 66671724
$ws.Range("K93").Value = 66671724
$ws.Range("M93").Value = -66670476

$ws.Range("H126").Value = 1990436.2
$ws.Range("I126").Value = 2651415
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 7954245
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -7951775
$ws.Range("N126").Value = -27440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 379395.22
$ws.Range("I122").Value = 518684.94
$ws.Range("J122").Value = 5054.0625
$ws.Range("K122").Value = 1556054.82
$ws.Range("L122").Value = 15162.1875
$ws.Range("M122").Value = -1553604.82
$ws.Range("N122").Value = -20062.1875
